# Generate Report for Handback
#
# For both the "zh-cn" and "de-de" worksheets, the two files that were
# pending handoff (a0c517ac-...md and f0e063cd-...md) have now been
# handed back and are in sync with en-US. This script:
#   - updates the Status column (B) from "Ready for handoff" to
#     "Handed back: in sync with en-US" for rows 2 and 3
#   - fills in the "Latest Target File" (E) and "Latest Handback File" (F)
#     columns with hyperlinked file names for rows 2 and 3
#   - stamps the "Latest Handback DateTime" (G) column for rows 2 and 3
#     with the handback timestamp

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

function Update-LangSheet($SheetName, $Md1Url, $Xlf1Url, $Md2Url, $Xlf2Url, $HandbackTime) {

    $ws = $wb.Worksheets.Item($SheetName)

    # Row 2 : a0c517ac-dd6c-45dc-8028-94bf10803aa3.md
    $ws.Range("B2").Value = $newStatus
    $ws.Hyperlinks.Add($ws.Range("E2"), $Md1Url, "", "", "a0c517ac-dd6c-45dc-8028-94bf10803aa3.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F2"), $Xlf1Url, "", "", "a0c517ac-dd6c-45dc-8028-94bf10803aa3.b060d66d737eec11c956d89fd8e5997b76f78c90.$SheetName.xlf") | Out-Null
    $ws.Range("G2").Value = $HandbackTime

    # Row 3 : f0e063cd-b80b-42ab-83e7-261f41b2a5aa.md
    $ws.Range("B3").Value = $newStatus
    $ws.Hyperlinks.Add($ws.Range("E3"), $Md2Url, "", "", "f0e063cd-b80b-42ab-83e7-261f41b2a5aa.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F3"), $Xlf2Url, "", "", "f0e063cd-b80b-42ab-83e7-261f41b2a5aa.20a172edba6dd3ea915b334ced012d9d9335dd8b.$SheetName.xlf") | Out-Null
    $ws.Range("G3").Value = $HandbackTime
}

# zh-cn sheet
Update-LangSheet "zh-cn" `
    "https://github.com/OpenLocalizationTest/oltest/blob/fef63f6d68fdc7df09dadf33d2fbcc32a9a98c80/e2e/a0c517ac-dd6c-45dc-8028-94bf10803aa3.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/90f45ea112180eb84def347972a2f3cefe7ce42d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a0c517ac-dd6c-45dc-8028-94bf10803aa3.b060d66d737eec11c956d89fd8e5997b76f78c90.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/fef63f6d68fdc7df09dadf33d2fbcc32a9a98c80/e2e/f0e063cd-b80b-42ab-83e7-261f41b2a5aa.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/90f45ea112180eb84def347972a2f3cefe7ce42d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f0e063cd-b80b-42ab-83e7-261f41b2a5aa.20a172edba6dd3ea915b334ced012d9d9335dd8b.zh-cn.xlf" `
    "2016-03-09 09:53:15"

# de-de sheet
Update-LangSheet "de-de" `
    "https://github.com/OpenLocalizationTest/oltest/blob/fef63f6d68fdc7df09dadf33d2fbcc32a9a98c80/e2e/a0c517ac-dd6c-45dc-8028-94bf10803aa3.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/71222215098bda3c67fae2a9c264c203bdec49ad/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a0c517ac-dd6c-45dc-8028-94bf10803aa3.b060d66d737eec11c956d89fd8e5997b76f78c90.de-de.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/fef63f6d68fdc7df09dadf33d2fbcc32a9a98c80/e2e/f0e063cd-b80b-42ab-83e7-261f41b2a5aa.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/71222215098bda3c67fae2a9c264c203bdec49ad/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f0e063cd-b80b-42ab-83e7-261f41b2a5aa.20a172edba6dd3ea915b334ced012d9d9335dd8b.de-de.xlf" `
    "2016-03-09 09:53:22"

Write-Output "Handback report generated."
